# Adjust spreader field efficiency (from 0.3 to 0.5) and spraying speed (from 15 to 20)
# on the "Mach 1" sheet (named ranges: sprayer_speed -> B101, spreader_eff -> B127)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mach 1")

# sprayer_speed: 15 -> 20
$ws.Range("B101").Value = 20

# spreader_eff: 0.3 -> 0.5
$ws.Range("B127").Value = 0.5
